$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("H2").Value = 40
$ws.Range("J3").Value = 63
$ws.Range("B9").Value = 132
$ws.Range("E9").Value = 156
$ws.Range("F9").Value = 179
$ws.Range("G9").Value = 182
$ws.Range("H9").Value = 139
$ws.Range("I9").Value = 188
$ws.Range("J9").Value = 138
$ws.Range("C10").Value = 440
$ws.Range("D10").Value = 577
$ws.Range("E10").Value = 644
$ws.Range("F10").Value = 751
$ws.Range("H10").Value = 152
$ws.Range("B11").Value = 524
$ws.Range("C11").Value = 646
$ws.Range("D11").Value = 808
$ws.Range("E11").Value = 878
$ws.Range("F11").Value = 999
$ws.Range("G11").Value = 748
$ws.Range("H11").Value = 368
$ws.Range("I11").Value = 562
$ws.Range("J11").Value = 498

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("F6").Value = 8
$ws.Range("F7").Value = 12

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("G6").Value = 20
$ws.Range("E7").Value = 23
$ws.Range("E8").Value = 39
$ws.Range("G8").Value = 46

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("H2").Value = 3
$ws.Range("J3").Value = 12
$ws.Range("H7").Value = 19
$ws.Range("I7").Value = 32
$ws.Range("C8").Value = 74
$ws.Range("D8").Value = 173
$ws.Range("E8").Value = 186
$ws.Range("F8").Value = 228
$ws.Range("C9").Value = 95
$ws.Range("D9").Value = 206
$ws.Range("E9").Value = 223
$ws.Range("F9").Value = 255
$ws.Range("H9").Value = 42
$ws.Range("I9").Value = 120
$ws.Range("J9").Value = 87

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("H2").Value = 1
$ws.Range("E6").Value = 18
$ws.Range("E7").Value = 26
$ws.Range("H7").Value = 10

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J5").Value = 3
$ws.Range("J6").Value = 5

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("E8").Value = 33
$ws.Range("C10").Value = 2
$ws.Range("E13").Value = 3
$ws.Range("F20").Value = 12
$ws.Range("J27").Value = 20
$ws.Range("E31").Value = 39
$ws.Range("G31").Value = 46
$ws.Range("F46").Value = 22
$ws.Range("I46").Value = 10
$ws.Range("F47").Value = 7
$ws.Range("J51").Value = 15
$ws.Range("C52").Value = 95
$ws.Range("D52").Value = 206
$ws.Range("E52").Value = 223
$ws.Range("F52").Value = 255
$ws.Range("H52").Value = 42
$ws.Range("I52").Value = 120
$ws.Range("J52").Value = 87
$ws.Range("E60").Value = 18
$ws.Range("J60").Value = 7
$ws.Range("B61").Value = 8
$ws.Range("E61").Value = 12
$ws.Range("G61").Value = 6
$ws.Range("J64").Value = 5
$ws.Range("E69").Value = 26
$ws.Range("H69").Value = 10
$ws.Range("E71").Value = 4
$ws.Range("E76").Value = 29
$ws.Range("H76").Value = 16
$ws.Range("D77").Value = 18
$ws.Range("E77").Value = 19
$ws.Range("C84").Value = 10
$ws.Range("C86").Value = 4
$ws.Range("B88").Value = 2
$ws.Range("E88").Value = 3
$ws.Range("B97").Value = 524
$ws.Range("C97").Value = 646
$ws.Range("D97").Value = 808
$ws.Range("E97").Value = 878
$ws.Range("F97").Value = 999
$ws.Range("G97").Value = 748
$ws.Range("H97").Value = 368
$ws.Range("I97").Value = 562
$ws.Range("J97").Value = 498

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("D5").Value = 16
$ws.Range("E5").Value = 18
$ws.Range("D6").Value = 18
$ws.Range("E6").Value = 19

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J4").Value = 6
$ws.Range("J6").Value = 20

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("F6").Value = 5
$ws.Range("I6").Value = 6
$ws.Range("F8").Value = 22
$ws.Range("I8").Value = 10

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 10
$ws.Range("J5").Value = 4

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("J7").Value = 15

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("B4").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("E5").Value = 8
$ws.Range("B6").Value = 8
$ws.Range("E6").Value = 12
$ws.Range("G6").Value = 6

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 3

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("E7").Value = 5
$ws.Range("H7").Value = 6
$ws.Range("H8").Value = 9
$ws.Range("E9").Value = 29
$ws.Range("H9").Value = 16

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("F4").Value = 1
$ws.Range("F6").Value = 7

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 4

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 2

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("B2").Value = 1
$ws.Range("E3").Value = 3
$ws.Range("B4").Value = 2
$ws.Range("E4").Value = 3

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("E6").Value = 16
$ws.Range("E7").Value = 33

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("C6").Value = 2
$ws.Range("C7").Value = 4
